$d = $word.ActiveDocument
$d.Content.Find.Execute("repetir la la instrucción de lectura.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "repetir la instrucción de lectura.", 2)
